$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.459.47'

$ws.Range('D3').Value = '1.960.06'
$ws.Range('E3').Value = '  -0.29%  '

$ws.Range('D4').Value = "'1.013"
$ws.Range('E4').Value = '  +0.41%  '

$ws.Range('D5').Value = "'322.01"
$ws.Range('E5').Value = '  -1.74%  '

$ws.Range('E6').Value = '  +0.17%  '

$ws.Range('D7').Value = "'0.4762"
$ws.Range('E7').Value = '  -4.56%  '

$ws.Range('D8').Value = "'0.4041"
$ws.Range('E8').Value = '  -4.29%  '

$ws.Range('D9').Value = "'53.97"
$ws.Range('E9').Value = '  +0.16%  '

$ws.Range('D10').Value = "'0.08443"
$ws.Range('E10').Value = '  -7.44%  '

$ws.Range('D11').Value = "'1.057"
$ws.Range('E11').Value = '  -3.86%  '

$ws.Range('D12').Value = "'22.31"
$ws.Range('E12').Value = '  -3.37%  '

$ws.Range('D13').Value = '2.001.37'
$ws.Range('E13').Value = '  +1.03%  '

$ws.Range('D14').Value = "'7.574"
$ws.Range('E14').Value = '  -4.12%  '

$ws.Range('D15').Value = "'6.163"
$ws.Range('E15').Value = '  -4.25%  '

$ws.Range('E16').Value = '  +0.30%  '

$ws.Range('D17').Value = "'90.72"
$ws.Range('E17').Value = '  -0.87%  '

$ws.Range('D18').Value = "'0.00001071"
$ws.Range('E18').Value = '  -2.89%  '

$ws.Range('D19').Value = "'0.06634"
$ws.Range('E19').Value = '  -0.29%  '

$ws.Range('D20').Value = "'18.52"
$ws.Range('E20').Value = '  -3.83%  '

$ws.Range('D21').Value = "'1.011"
$ws.Range('E21').Value = '  +0.22%  '

$ws.Range('D22').Value = "'5.851"
$ws.Range('E22').Value = '  -1.57%  '

$ws.Range('D23').Value = '28.520.33'
$ws.Range('E23').Value = '  -2.17%  '

$ws.Range('D24').Value = "'11.44"
$ws.Range('E24').Value = '  -4.70%  '

$ws.Range('D25').Value = "'2.299"
$ws.Range('E25').Value = '  +0.27%  '

$ws.Range('D26').Value = '2.229.57'
$ws.Range('E26').Value = '  +0.76%  '

$ws.Range('D27').Value = "'155.57"
$ws.Range('E27').Value = '  -0.97%  '

$ws.Range('D28').Value = "'20.26"
$ws.Range('E28').Value = '  -1.92%  '

$ws.Range('D29').Value = "'5.884"
$ws.Range('E29').Value = '  -5.80%  '

$ws.Range('D30').Value = "'2.158"
$ws.Range('E30').Value = '  -5.06%  '

$ws.Range('D31').Value = "'124.46"
$ws.Range('E31').Value = '  -2.31%  '

$ws.Range('D32').Value = "'0.9794"
$ws.Range('E32').Value = '  -6.13%  '

$ws.Range('D33').Value = "'0.09645"
$ws.Range('E33').Value = '  -2.13%  '

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = "'3.695"
$ws.Range('E34').Value = '  -0.37%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'1.446"
$ws.Range('E35').Value = '  -5.39%  '

$ws.Range('D36').Value = "'5.610"
$ws.Range('E36').Value = '  -3.69%  '

$ws.Range('D37').Value = "'9.151"
$ws.Range('E37').Value = '  +1.34%  '

$ws.Range('D38').Value = "'0.02328"
$ws.Range('E38').Value = '  -4.27%  '

$ws.Range('D39').Value = "'0.06220"
$ws.Range('E39').Value = '  -2.17%  '

$ws.Range('D40').Value = "'1.252"
$ws.Range('E40').Value = '  -2.80%  '

$ws.Range('D41').Value = "'0.6199"
$ws.Range('E41').Value = '  -3.87%  '

$ws.Range('D42').Value = "'11.14"
$ws.Range('E42').Value = '  -3.14%  '

$ws.Range('D43').Value = "'1.011"
$ws.Range('E43').Value = '  +0.21%  '

$ws.Range('D44').Value = "'0.1909"
$ws.Range('E44').Value = '  -4.83%  '

$ws.Range('D45').Value = "'1.355"
$ws.Range('E45').Value = '  +5.99%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.5935"
$ws.Range('E46').Value = '  -4.72%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'12.98"
$ws.Range('E47').Value = '  -3.49%  '

$ws.Range('D48').Value = "'2.050"
$ws.Range('E48').Value = '  -6.09%  '

$ws.Range('D49').Value = "'3.395"
$ws.Range('E49').Value = '  -2.32%  '

$ws.Range('D50').Value = "'0.06816"
$ws.Range('E50').Value = '  -0.75%  '

$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = "'111.15"
$ws.Range('E51').Value = '  -1.47%  '
